$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0314455509185791
$ws.Range("C2").Value = 0.04860405921936035
$ws.Range("D2").Value = 0.01070356369018555
$ws.Range("E2").Value = 0.03422327041625976
$ws.Range("F2").Value = 0.002495002746582031
$ws.Range("G2").Value = 0.1478825092315674
$ws.Range("H2").Value = 0.03303136825561524
$ws.Range("I2").Value = 0.06153683662414551
$ws.Range("J2").Value = 0.0169736385345459
$ws.Range("K2").Value = 0.05254664421081543
$ws.Range("L2").Value = 0.004017877578735352
$ws.Range("M2").Value = 0.04060540199279785
$ws.Range("B3").Value = 0.147467565536499
$ws.Range("C3").Value = 0.04843978881835938
$ws.Range("D3").Value = 0.0328115463256836
$ws.Range("E3").Value = 0.0176877498626709
$ws.Range("F3").Value = 0.01966137886047363
$ws.Range("G3").Value = 0.01844062805175781
$ws.Range("H3").Value = 0.1804601669311524
$ws.Range("I3").Value = 0.0492513656616211
$ws.Range("J3").Value = 0.1045034885406494
$ws.Range("K3").Value = 0.03377041816711426
$ws.Range("L3").Value = 0.03438029289245605
$ws.Range("M3").Value = 0.02857403755187988
$ws.Range("B4").Value = 0.05333418846130371
$ws.Range("C4").Value = 0.02750449180603027
$ws.Range("D4").Value = 0.01539821624755859
$ws.Range("E4").Value = 0.01601839065551758
$ws.Range("F4").Value = 0.06522397994995117
$ws.Range("G4").Value = 0.01664719581604004
$ws.Range("H4").Value = 0.03271093368530274
$ws.Range("I4").Value = 0.03025169372558594
$ws.Range("J4").Value = 0.02226753234863281
$ws.Range("K4").Value = 0.02833094596862793
$ws.Range("L4").Value = 0.07210717201232911
$ws.Range("M4").Value = 0.01884493827819824
$ws.Range("B5").Value = 0.02825822830200195
$ws.Range("C5").Value = 0.02827963829040527
$ws.Range("D5").Value = 0.01311135292053223
$ws.Range("E5").Value = 0.02474184036254883
$ws.Range("H5").Value = 0.02676587104797363
$ws.Range("I5").Value = 0.02411022186279297
$ws.Range("J5").Value = 0.01754975318908691
$ws.Range("K5").Value = 0.02553768157958984
$ws.Range("B6").Value = 0.9008755207061767
$ws.Range("C6").Value = 0.04585442543029785
$ws.Range("D6").Value = 0.8196640968322754
$ws.Range("E6").Value = 0.05265550613403321
$ws.Range("F6").Value = 0.9259534358978272
$ws.Range("G6").Value = 0.0311607837677002
$ws.Range("H6").Value = 0.2254391193389892
$ws.Range("I6").Value = 0.02995009422302246
$ws.Range("J6").Value = 0.278513765335083
$ws.Range("K6").Value = 0.03061251640319824
$ws.Range("L6").Value = 0.9483530521392822
$ws.Range("M6").Value = 0.031549072265625
$ws.Range("B7").Value = 1.357732725143433
$ws.Range("C7").Value = 0.1557707786560059
$ws.Range("D7").Value = 0.5387883663177491
$ws.Range("E7").Value = 0.05468273162841797
$ws.Range("F7").Value = 0.7940410614013672
$ws.Range("G7").Value = 0.04458446502685547
$ws.Range("H7").Value = 1.066512250900268
$ws.Range("I7").Value = 0.1071974754333496
$ws.Range("J7").Value = 0.5797893524169921
$ws.Range("K7").Value = 0.05726819038391114
$ws.Range("L7").Value = 0.9662599086761474
$ws.Range("M7").Value = 0.04015698432922363